# TC_63812 - Updated base change test case data
# Applies the OOXML diff changes to the "Add Devices Loop A" sheet:
#  - Adds a new "IBUnitsLoadingDetail" column (I)
#  - Converts several label cells (F8, F9, D10, F10, D11) from descriptive
#    text to plain numeric "row index" values
#  - Preserves the old descriptive text as cell comments (notes)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use the same author name as the original commit for any comments added.
$excel.UserName = "Alpesh Dhakad"

# --- New column header (I7) -------------------------------------------------
# Copy formatting from the neighbouring header cell (H7) so the new header
# cell picks up the same font/fill/border style, then set its text.
$ws.Cells.Item(7, 8).Copy($ws.Cells.Item(7, 9))
$ws.Cells.Item(7, 9).Value = "IBUnitsLoadingDetail"

# --- Row 8 (460PH) -----------------------------------------------------------
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(8, 9))
$ws.Cells.Item(8, 9).Value = "460PH - 1"
$ws.Cells.Item(8, 9).Value = "460PH - 1"

# --- Row 9 (460P) ------------------------------------------------------------
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(9, 9))
$ws.Cells.Item(9, 9).Value = "460PH - 1"

# --- Row 10 (460H) -----------------------------------------------------------
$ws.Cells.Item(10, 4).Value = 28
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(10, 9))
$ws.Cells.Item(10, 9).Value = "460PH - 1"

# --- Row 11 (460PC) ----------------------------------------------------------
$ws.Cells.Item(11, 4).Value = 16
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(11, 9))
$ws.Cells.Item(11, 9).Value = "460PH - 1"

# --- Row 12 (410RIM) ---------------------------------------------------------
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(12, 9))
$ws.Cells.Item(12, 9).Value = "460PH - 1"

# --- Row 13 (420CP) ----------------------------------------------------------
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(13, 9))
$ws.Cells.Item(13, 9).Value = "460PH - 1"

# --- Row 14 (410MIM) ---------------------------------------------------------
$ws.Cells.Item(7, 7).Copy($ws.Cells.Item(14, 9))
$ws.Cells.Item(14, 9).Value = "460PH - 1"

# --- Comments preserving the previous descriptive values --------------------
$ws.Range("F8").AddComment("Alpesh Dhakad:" + [char]10 + "4B-EM 4`" [517.050.052]") | Out-Null
$ws.Range("F9").AddComment("Alpesh Dhakad:" + [char]10 + "4B-EM 4`" [517.050.052]") | Out-Null
$ws.Range("D10").AddComment("Alpesh Dhakad:" + [char]10 + "440DSB [576.440.001] @ 60dB") | Out-Null
$ws.Range("F10").AddComment("Alpesh Dhakad:" + [char]10 + "A-CON [557.080.002]") | Out-Null
$ws.Range("D11").AddComment("Alpesh Dhakad:" + [char]10 + "430SB [516.800.710] @ 68dB") | Out-Null

# --- Column width for the new column -----------------------------------------
$ws.Columns.Item(9).AutoFit()

# --- Selection matches the state left behind in the source workbook ---------
$ws.Range("I12").Select()
